$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 839.9
$ws.Range("I18").Value = 839.9
$ws.Range("K18").Value = 839.9
$ws.Range("M18").Value = -555.9
# Row 97
$ws.Range("H97").Value = 3890.9092
$ws.Range("J97").Value = 3890.9092
$ws.Range("L97").Value = 11672.7276
$ws.Range("N97").Value = -12664.7276
# Row 127
$ws.Range("H127").Value = 3138.5715
$ws.Range("I127").Value = 3138.5715
$ws.Range("K127").Value = 9415.7145
$ws.Range("M127").Value = -4455.7145
# Row 138
$ws.Range("H138").Value = 2417.7778
$ws.Range("I138").Value = 1339.48
$ws.Range("J138").Value = 2991.3403
$ws.Range("K138").Value = 4018.44
$ws.Range("L138").Value = 8974.0209
$ws.Range("M138").Value = 1121.56
$ws.Range("N138").Value = -19254.0209

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9262828
$ws.Range("I32").Value = 10206973
$ws.Range("J32").Value = 10202
$ws.Range("K32").Value = 10206973
$ws.Range("L32").Value = 10202
$ws.Range("M32").Value = -10206686
$ws.Range("N32").Value = -10776
# Row 61
$ws.Range("H61").Value = 10443757
$ws.Range("I61").Value = 14708917
$ws.Range("K61").Value = 14708917
$ws.Range("M61").Value = -14708705
# Row 63
$ws.Range("H63").Value = 4015.4
$ws.Range("I63").Value = 4127.5557
$ws.Range("J63").Value = 3006
$ws.Range("K63").Value = 4127.5557
$ws.Range("L63").Value = 3006
$ws.Range("M63").Value = -3441.5557
$ws.Range("N63").Value = -4378
# Row 66
$ws.Range("H66").Value = 4015.4
$ws.Range("I66").Value = 4127.5557
$ws.Range("J66").Value = 3006
$ws.Range("K66").Value = 20637.7785
$ws.Range("L66").Value = 15030
$ws.Range("M66").Value = -17205.7785
$ws.Range("N66").Value = -21894
# Row 74
$ws.Range("H74").Value = 5107398.5
$ws.Range("I74").Value = 6251346
$ws.Range("K74").Value = 6251346
$ws.Range("M74").Value = -6250472
# Row 77
$ws.Range("H77").Value = 5107398.5
$ws.Range("I77").Value = 6251346
$ws.Range("K77").Value = 31256730
$ws.Range("M77").Value = -31252362
# Row 88
$ws.Range("H88").Value = 1885.6
$ws.Range("I88").Value = 2197
$ws.Range("J88").Value = 1772.3636
$ws.Range("K88").Value = 2197
$ws.Range("L88").Value = 1772.3636
$ws.Range("M88").Value = -1791
$ws.Range("N88").Value = -2584.3636
# Row 91
$ws.Range("H91").Value = 1885.6
$ws.Range("I91").Value = 2197
$ws.Range("J91").Value = 1772.3636
$ws.Range("K91").Value = 2197
$ws.Range("L91").Value = 1772.3636
$ws.Range("M91").Value = -793
$ws.Range("N91").Value = -4580.3636
# Row 132
$ws.Range("H132").Value = 3410.5557
$ws.Range("I132").Value = 1436.2
$ws.Range("J132").Value = 5878.5
$ws.Range("K132").Value = 4308.6
$ws.Range("L132").Value = 17635.5
$ws.Range("M132").Value = -1778.6
$ws.Range("N132").Value = -22695.5
# Row 136
$ws.Range("H136").Value = 10443757
$ws.Range("I136").Value = 14708917
$ws.Range("K136").Value = 44126751
$ws.Range("M136").Value = -44124201

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2241.2
$ws.Range("I105").Value = 1585.25
$ws.Range("K105").Value = 1585.25
$ws.Range("M105").Value = 161.75
# Row 134
$ws.Range("H134").Value = 27660.205
$ws.Range("I134").Value = 1621.9395
$ws.Range("J134").Value = 170870.67
$ws.Range("K134").Value = 4865.818499999999
$ws.Range("L134").Value = 512612.01
$ws.Range("M134").Value = -2330.818499999999
$ws.Range("N134").Value = -517682.01

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 537858.75
$ws.Range("I31").Value = 10702.05
$ws.Range("J31").Value = 977155.9399999999
$ws.Range("K31").Value = 10702.05
$ws.Range("L31").Value = 977155.9399999999
$ws.Range("M31").Value = -10407.05
$ws.Range("N31").Value = -977745.9399999999
# Row 34
$ws.Range("H34").Value = 537858.75
$ws.Range("I34").Value = 10702.05
$ws.Range("J34").Value = 977155.9399999999
$ws.Range("K34").Value = 10702.05
$ws.Range("L34").Value = 977155.9399999999
$ws.Range("M34").Value = -10500.05
$ws.Range("N34").Value = -977559.9399999999
# Row 94
$ws.Range("H94").Value = 5781.727
$ws.Range("I94").Value = 5872.6665
$ws.Range("J94").Value = 5747.625
$ws.Range("K94").Value = 5872.6665
$ws.Range("L94").Value = 5747.625
$ws.Range("M94").Value = -5421.6665
$ws.Range("N94").Value = -6649.625
# Row 110
$ws.Range("H110").Value = 109964.5
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 109964.5
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 109964.5
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -118144.5
# Row 122
$ws.Range("H122").Value = 1545.4286
$ws.Range("I122").Value = 1620.8
$ws.Range("J122").Value = 1357
$ws.Range("K122").Value = 4862.4
$ws.Range("L122").Value = 4071
$ws.Range("M122").Value = -2412.4
$ws.Range("N122").Value = -8971
# Row 132
$ws.Range("H132").Value = 2908.1936
$ws.Range("I132").Value = 2221.1738
$ws.Range("J132").Value = 4883.375
$ws.Range("K132").Value = 6663.5214
$ws.Range("L132").Value = 14650.125
$ws.Range("M132").Value = -4133.5214
$ws.Range("N132").Value = -19710.125
# Row 134
$ws.Range("H134").Value = 436347.2
$ws.Range("I134").Value = 501517.66
$ws.Range("J134").Value = 1877.3334
$ws.Range("K134").Value = 1504552.98
$ws.Range("L134").Value = 5632.0002
$ws.Range("M134").Value = -1502017.98
$ws.Range("N134").Value = -10702.0002

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 5479.8
$ws.Range("I80").Value = 5449.5
$ws.Range("K80").Value = 16348.5
$ws.Range("M80").Value = -15412.5
# Row 83
$ws.Range("H83").Value = 5479.8
$ws.Range("I83").Value = 5449.5
$ws.Range("K83").Value = 49045.5
$ws.Range("M83").Value = -44365.5
# Row 107
$ws.Range("H107").Value = 575.2
$ws.Range("J107").Value = 662.4
$ws.Range("L107").Value = 1987.2
$ws.Range("N107").Value = -5827.2

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 109
$ws.Range("H109").Value = 37642.5
$ws.Range("J109").Value = 37642.5
$ws.Range("L109").Value = 37642.5
$ws.Range("N109").Value = -39722.5
# Row 132
$ws.Range("H132").Value = 31254580
$ws.Range("I132").Value = 35715812
$ws.Range("K132").Value = 107147436
$ws.Range("M132").Value = -107144906

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 56767.316
$ws.Range("I7").Value = 3288.5833
$ws.Range("K7").Value = 3288.5833
$ws.Range("M7").Value = -3176.5833
# Row 126
$ws.Range("H126").Value = 56767.316
$ws.Range("I126").Value = 3288.5833
$ws.Range("K126").Value = 9865.749899999999
$ws.Range("M126").Value = -7395.749899999999
# Row 132
$ws.Range("H132").Value = 23474.576
$ws.Range("I132").Value = 5082.1943
$ws.Range("K132").Value = 15246.5829
$ws.Range("M132").Value = -12716.5829
# Row 136
$ws.Range("H136").Value = 57502.74
$ws.Range("I136").Value = 10597.091
$ws.Range("K136").Value = 31791.273
$ws.Range("M136").Value = -29241.273

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2513.6365
$ws.Range("I132").Value = 2065
$ws.Range("K132").Value = 6195
$ws.Range("M132").Value = -3665
# Row 136
$ws.Range("H136").Value = 2319
$ws.Range("I136").Value = 726
$ws.Range("J136").Value = 6415.2856
$ws.Range("K136").Value = 2178
$ws.Range("L136").Value = 19245.8568
$ws.Range("M136").Value = 372
$ws.Range("N136").Value = -24345.8568
